$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.571.88'
$ws.Range("E2").Value = '  -1.11%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.444.55'
$ws.Range("E3").Value = '  -3.47%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '593.48'
$ws.Range("E5").Value = '  -1.93%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '136.20'
$ws.Range("E6").Value = '  -7.65%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.442.88'
$ws.Range("E7").Value = '  -3.45%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.11%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.492'
$ws.Range("E9").Value = '  +0.40%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.40'
$ws.Range("E10").Value = '  -6.30%  '
$ws.Range("E11").Value = '  -9.26%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.379'
$ws.Range("E12").Value = '  -7.59%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.025.65'
$ws.Range("E13").Value = '  -3.59%  '
$ws.Range("E14").Value = '  -10.68%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '26.67'
$ws.Range("E15").Value = '  -9.29%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.458.51'
$ws.Range("E16").Value = '  -3.22%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '65.556.19'
$ws.Range("E17").Value = '  -1.11%  '
$ws.Range("E18").Value = '  -2.09%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.91'
$ws.Range("E19").Value = '  -9.66%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.82'
$ws.Range("E20").Value = '  -7.49%  '
$ws.Range("E21").Value = '  -7.02%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '395.13'
$ws.Range("E22").Value = '  -6.16%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.547'
$ws.Range("E23").Value = '  -10.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '73.32'
$ws.Range("E24").Value = '  -5.81%  '
$ws.Range("E25").Value = '  -0.09%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.587.11'
$ws.Range("E26").Value = '  -3.25%  '
$ws.Range("E27").Value = '  -10.08%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.998'
$ws.Range("E28").Value = '  -0.18%  '
$ws.Range("E29").Value = '  -8.73%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.20'
$ws.Range("E30").Value = '  -10.01%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.20'
$ws.Range("E31").Value = '  -11.77%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.450.03'
$ws.Range("E32").Value = '  -3.23%  '
$ws.Range("E33").Value = '  -0.03%  '
$ws.Range("E34").Value = '  -6.44%  '
$ws.Range("E35").Value = '  -7.32%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '172.46'
$ws.Range("E36").Value = '  -1.58%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.21'
$ws.Range("E37").Value = '  -13.16%  '
$ws.Range("E38").Value = '  -10.26%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.51'
$ws.Range("E39").Value = '  -7.37%  '
$ws.Range("E40").Value = '  -11.40%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0774'
$ws.Range("E41").Value = '  -8.10%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.823'
$ws.Range("E42").Value = '  -5.84%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '43.45'
$ws.Range("E43").Value = '  -5.14%  '
$ws.Range("E44").Value = '  +0.24%  '
$ws.Range("E45").Value = '  -14.24%  '
$ws.Range("E46").Value = '  -12.14%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '22.93'
$ws.Range("E47").Value = '  -2.29%  '
$ws.Range("E48").Value = '  -2.19%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '6.58'
$ws.Range("E49").Value = '  -7.34%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.10'
$ws.Range("E50").Value = '  -15.01%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.201.53'
$ws.Range("E51").Value = '  -7.74%  '
